# The "TECHNICAL SKILLS" table (category / spacer / skill columns, with
# {category} and {skill} merge-field placeholders) should render with a
# smaller font (9.5pt / half-point value 19) so the generated resume is
# lighter / more compact. Apply the size to both the Latin (Size) and
# complex-script (SizeBi) font properties for every paragraph mark and run
# in the table, which is what produces the <w:sz>/<w:szCs> pair on each
# run/paragraph properties element.

$d = $word.ActiveDocument

$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text -like "*{category}*" -and $candidate.Range.Text -like "*{skill}*") {
        $targetTable = $candidate
        break
    }
}

if ($targetTable -ne $null) {
    $targetTable.Range.Font.Size = 9.5
    $targetTable.Range.Font.SizeBi = 9.5
}
